# Update mass-flow results with newly re-run input-file values.
# Output_flows and Input_flows sheets get new simulation outputs for a
# subset of particle-size / behaviour-category rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Output_flows")
$ws.Range("C2").Value = 0.000000000000000465047653060642
$ws.Range("E2").Value = 0.0000000000002381623467373324
$ws.Range("G2").Value = 0.00000000000001118203238079399
$ws.Range("I2").Value = 0.00000000000000004136526898892864
$ws.Range("K2").Value = 0.0000000000000000006228524882861391
$ws.Range("C7").Value = 0.00000004817938461539514
$ws.Range("F7").Value = 0.0000000000002498514148926642
$ws.Range("G7").Value = 0.00001158469320961876
$ws.Range("I7").Value = 0.00000004285481694661923
$ws.Range("J7").Value = 0.00000001028310872941287
$ws.Range("C12").Value = 0.00000000003315434224269436
$ws.Range("E12").Value = 0.00003764720642638757
$ws.Range("I12").Value = 0.00000000001474513714878609
$ws.Range("J12").Value = 0.0000000007666735092058083
$ws.Range("C13").Value = 0.00000000000376063442136776
$ws.Range("D13").Value = 0.00000000000278704552926555
$ws.Range("E13").Value = 0.00001592603812282098
$ws.Range("I13").Value = 0.000000000001672513057372816
$ws.Range("J13").Value = 0.0000000001700823802208533
$ws.Range("C14").Value = 0.0000000000004487038797241835
$ws.Range("D14").Value = 0.000000000001330156565970652
$ws.Range("E14").Value = 0.0001629826987669894
$ws.Range("I14").Value = 0.0000000000001995575782289392
$ws.Range("J14").Value = 0.0000000001177158598848538
$ws.Range("C17").Value = 0.008286139366879166
$ws.Range("F17").Value = 0.00003764727074361214
$ws.Range("I17").Value = 0.01474078544574493
$ws.Range("J17").Value = 0.7666424305967253
$ws.Range("C18").Value = 0.0009400568567891301
$ws.Range("D18").Value = 0.00002786778305416879
$ws.Range("F18").Value = 0.00001592605507075637
$ws.Range("I18").Value = 0.001672332049846877
$ws.Range("J18").Value = 0.1700795654137782
$ws.Range("C19").Value = 0.0001121754757839671
$ws.Range("D19").Value = 0.00001330150963077166
$ws.Range("F19").Value = 0.0001629827102093696
$ws.Range("I19").Value = 0.0001995566991565819
$ws.Range("J19").Value = 0.1177158464177236

$ws = $wb.Worksheets.Item("Input_flows")
$ws.Range("C7").Value = 0.00001168541664970464
$ws.Range("C12").Value = 0.0000000007230539365951427
$ws.Range("C13").Value = 0.0000000001555650373312833
$ws.Range("C14").Value = 0.0000000001055062157244061
$ws.Range("C17").Value = 0.7070113045943595
$ws.Range("C18").Value = 0.1543930611725542
$ws.Range("C19").Value = 0.1055014646903124
$ws.Range("E27").Value = 0.0000000005938818946236801
$ws.Range("E32").Value = 0.00000000002373852157806095
$ws.Range("E33").Value = 0.00000000000512452223141207
$ws.Range("E34").Value = 0.000000000002745682059389261
$ws.Range("E37").Value = 0.08261194764822795
$ws.Range("E38").Value = 0.01832011019304685
$ws.Range("E39").Value = 0.01253941542342499
